$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SwateTemplateMetadata")

# Rename the metadata sheet to isa_template
$ws.Name = "isa_template"

# Remove the GEO-specific ER/Tag values from the generic isa_template sheet
# (row 12: Tags, row 13: Tags Term Accession Number, row 14: Tags Term Source REF)
# by shifting each row's contents one column to the left, dropping column B.

# Remove hyperlink currently anchored at G13 before we shift cell contents,
# then re-add it at its new location (F13) once the shift is done.
$ws.Range("G13").Hyperlinks.Delete()

# Row 12 has data through column H -> shift C:H into B:G
$ws.Cells.Item(12,3).Cut($ws.Cells.Item(12,2))
$ws.Cells.Item(12,4).Cut($ws.Cells.Item(12,3))
$ws.Cells.Item(12,5).Cut($ws.Cells.Item(12,4))
$ws.Cells.Item(12,6).Cut($ws.Cells.Item(12,5))
$ws.Cells.Item(12,7).Cut($ws.Cells.Item(12,6))
$ws.Cells.Item(12,8).Cut($ws.Cells.Item(12,7))
$ws.Cells.Item(12,8).Clear()

# Row 13 has data through column G -> shift C:G into B:F
$ws.Cells.Item(13,3).Cut($ws.Cells.Item(13,2))
$ws.Cells.Item(13,4).Cut($ws.Cells.Item(13,3))
$ws.Cells.Item(13,5).Cut($ws.Cells.Item(13,4))
$ws.Cells.Item(13,6).Cut($ws.Cells.Item(13,5))
$ws.Cells.Item(13,7).Cut($ws.Cells.Item(13,6))
$ws.Cells.Item(13,7).Clear()

# Row 14 has data through column G -> shift C:G into B:F
$ws.Cells.Item(14,3).Cut($ws.Cells.Item(14,2))
$ws.Cells.Item(14,4).Cut($ws.Cells.Item(14,3))
$ws.Cells.Item(14,5).Cut($ws.Cells.Item(14,4))
$ws.Cells.Item(14,6).Cut($ws.Cells.Item(14,5))
$ws.Cells.Item(14,7).Cut($ws.Cells.Item(14,6))
$ws.Cells.Item(14,7).Clear()

# Re-create the hyperlink at its new home, F13
$ws.Hyperlinks.Add($ws.Range("F13"), "http://purl.obolibrary.org/obo/NCIT_C153189") | Out-Null

# Restore selection/active cell to B9 on the renamed sheet
$ws.Activate()
$ws.Range("B9").Select()
